$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("issues")

# Update the "Assignee" column (D) from single-name values to
# comma-separated multi-assignee lists. Statement order below is chosen
# so that the new shared-string entries are first created in the same
# order they appear in the target workbook (Richard's group, then
# Daniel's group, then Danny's group), matching the append order Excel
# uses when rebuilding the shared strings table on save.

$richard = "RichardWaiteSTFC, martyngigg, SilkeSchomann, sf1919, ConorMFinn"
$daniel  = "DanielMurphy22, gemmaguest, AnthonyLim23, Pasarus"
$danny   = "DannyHindson, MialLewis, DavidFair, Harrietbrown, thomashampson"

# First usages, in the order they should be introduced into the shared
# string table.
$ws.Range("D3").Value = $richard
$ws.Range("D5").Value = $daniel
$ws.Range("D2").Value = $danny

# Remaining rows that reuse the same three strings.
$ws.Range("D6").Value = $richard
$ws.Range("D7").Value = $richard
$ws.Range("D8").Value = $richard
$ws.Range("D15").Value = $richard
$ws.Range("D18").Value = $richard

$ws.Range("D11").Value = $daniel
$ws.Range("D13").Value = $daniel
$ws.Range("D14").Value = $daniel
$ws.Range("D16").Value = $daniel
$ws.Range("D17").Value = $daniel

$ws.Range("D4").Value = $danny
$ws.Range("D9").Value = $danny
$ws.Range("D10").Value = $danny
$ws.Range("D12").Value = $danny

# Rename the "Manual Testing ISIS SANS (new GUI)" test to drop the
# "(new GUI)" suffix. This is the last new shared string introduced.
$ws.Range("A12").Value = "Manual Testing ISIS SANS"

# The longer assignee lists no longer fit on one line at the existing
# row heights, so re-apply the row heights Excel would compute when
# auto-fitting wrapped text in column D for these rows.
$ws.Rows.Item(4).RowHeight = 48
$ws.Rows.Item(9).RowHeight = 48
$ws.Rows.Item(11).RowHeight = 48
$ws.Rows.Item(12).RowHeight = 48
$ws.Rows.Item(13).RowHeight = 48
$ws.Rows.Item(16).RowHeight = 48
$ws.Rows.Item(17).RowHeight = 48

# Update the current selection/view to match the saved workbook state.
$ws.Activate()
$ws.Range("A13").Select()
